# "Common: Added boring stuff for liquid"
# Append six new vendor names to the "vendors" sheet and re-sort the
# vendor list (A2:A81) alphabetically, exactly like a user typing the
# new rows at the bottom and then re-applying Data > Sort.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

# New vendors, appended in the order they were typed (this also controls
# the order in which they land in the shared-strings table).
$newVendors = @(
    'Tobacco Bastards',
    'Rocket Girl',
    'Electra',
    'Five Pawns',
    'MaZa',
    'Ripe Vapes'
)

$firstNewRow = 76
for ($i = 0; $i -lt $newVendors.Length; $i++) {
    $ws.Cells.Item($firstNewRow + $i, 1).Value = $newVendors[$i]
}

$lastRow = $firstNewRow + $newVendors.Length - 1

# Re-sort A2:A81 (column A, ascending, header in row 1) -- this reorders
# the rows into alphabetical order just like the original sort did.
$sortRange = $ws.Range("A2:A" + $lastRow)
$sortRange.Sort($ws.Range("A2")) | Out-Null

# Keep the worksheet's recorded sort state/range in sync with the new
# extent of the list.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A" + $lastRow)) | Out-Null
$sortObj.SetRange($ws.Range("A2:A" + $lastRow))
$sortObj.Header = 2
$sortObj.Apply()

# Match the resulting selection/scroll position left behind in the file.
$ws.Range("A56").Select() | Out-Null

Write-Output ("vendors sheet now has {0} rows (A2:A{1})" -f $newVendors.Length, $lastRow)
